$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(11961, "18PG",  0, 9.3028, 3, 10, 36, 1, 10, 70,  75, 180, -0.005874720885638602, 1.946117060036923),
    @(305,   "18PG",  0, 9.3028, 3, 10, 36, 1, 5,  350, 45, 4,   -0.01045724615189837,  1.904086138294788),
    @(9159,  "18PG",  0, 9.3028, 3, 10, 36, 1, 9,  50,  65, 120, -0.0005510610013151772, 1.994945668495937),
    @(21292, "DOPE",  1, 9.784000000000001, 2, 9, 36, 2, 7, 400, 75, 140, -0.03151380054466597, 1.710955421404324),
    @(21125, "DOPE",  1, 9.784000000000001, 2, 9, 36, 2, 7, 450, 55, 120, -0.00662968209486068, 1.939192555825938),
    @(33632, "DOTAP", 1, 12.515, 0, 4, 36, 2, 7,  60,  75, 2, 0.8287592435514013, 9.601379781853453),
    @(28991, "DOTAP", 1, 12.515, 0, 4, 36, 2, 5,  350, 25, 4, 0.7933505544661379, 9.276611285563417),
    @(33974, "DOTAP", 1, 12.515, 0, 4, 36, 2, 7,  250, 55, 5, 0.7423982913727166, 8.809277128470558),
    @(33663, "DOTAP", 1, 12.515, 0, 4, 36, 2, 7,  90,  30, 3, 0.6965966909481704, 8.389184849376621),
    @(28874, "DOTAP", 1, 12.515, 0, 4, 36, 2, 5,  80,  30, 3, 0.7517572385919185, 8.895117392365076),
    @(38313, "DOTAP", 1, 12.515, 0, 4, 36, 2, 10, 200, 25, 2, 0.6965966909481704, 8.389184849376621),
    @(43165, "DSPC",  1, 6.12,   0, 8, 36, 0, 5,  50,  55, 2, -0.03976853977306633, 1.635242953201436)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $ws.Cells.Item($row, 11).Value = $r[10]
    $ws.Cells.Item($row, 12).Value = $r[11]
    $ws.Cells.Item($row, 13).Value = $r[12]
    $ws.Cells.Item($row, 14).Value = $r[13]
    $row++
}
